# 10.02/2024 - return to host 21
# This script applies the updates described in the commit:
#  - Column I (Date_of_sales) for rows 2..30: 45341 -> 45342 (advance by one day)
#  - Column C (short Param text) reorders some tokens for several rows
#  - Column G (comma separated Param list) reorders the same tokens for several rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

# --- Column I: bump every date in the data rows (2-30) by one day ---
for ($row = 2; $row -le 30; $row++) {
    $cell = $ws.Cells.Item($row, 9)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        $cell.Value = $cur + 1
    }
}

# --- Column C: reordered short param strings ---
$colC = @{
    3  = "б/к легк сер"
    4  = "б/к легк сер"
    11 = "Type LS-2 груз сер"
    13 = "202B H Type LS-2 C сер"
    14 = "б/к груз сер"
    15 = "б/к легк сер"
    16 = "б/к легк сер"
    17 = "8 сер сх"
}
foreach ($row in $colC.Keys) {
    $ws.Cells.Item($row, 3).Value = $colC[$row]
}

# --- Column G: reordered comma-separated param lists ---
$colG = @{
    4  = "б/к, легк, сер"
    5  = "б/к, легк, сер"
    13 = "Type, LS-2, груз, сер"
    15 = "202B, H, Type, LS-2, C, сер"
    16 = "202B, H, Type, LS-2, C, сер"
    17 = "202B, H, Type, LS-2, C, сер"
    18 = "б/к, груз, сер"
    19 = "б/к, груз, сер"
    20 = "б/к, груз, сер"
    21 = "б/к, груз, сер"
    22 = "б/к, легк, сер"
    23 = "б/к, легк, сер"
}
foreach ($row in $colG.Keys) {
    $ws.Cells.Item($row, 7).Value = $colG[$row]
}
